$wb = $excel.ActiveWorkbook

# --- 1. Fix the "Decomposition" sheet view: reset the scrolled-down viewport ---
$wsDecomp = $wb.Worksheets.Item("Decomposition")
$wsDecomp.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$wsDecomp.Range("C24:E24").Select() | Out-Null

# --- 2. Update "Sheet1": insert a new title row, re-style, merge, and resize rows ---
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Insert a new row 1 above the existing header row, pushing everything down.
$ws.Rows.Item(1).Insert()

# Add merged title in the new row 1 spanning B1:E1.
$ws.Range("B1:E1").Merge()
$ws.Range("B1").Value = "Imact on PV ERC, plan starting with 100% funded ratio"
$ws.Range("B1:E1").HorizontalAlignment = -4108
$ws.Rows.Item(1).RowHeight = 28.5

# Bold the column-header row (now row 2: "Total effect" / "cost effect" / "risk effect").
$ws.Range("C2:E2").Font.Bold = $true

# Give each data row (now rows 3-10) a taller custom height.
$ws.Range("B3:E10").RowHeight = 24

# Update the selection to match the target state.
$ws.Range("B1:E10").Select()

$wb.Worksheets.Item("Sheet1").Activate()
